$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("MPSP")
$ws.Range("C2").Value = 2.622343602872005
$ws.Range("D2").Value = 2.152372601055987
$ws.Range("E2").Value = 1.792812421528958
$ws.Range("F2").Value = 1.369815210928913
$ws.Range("C3").Value = 2.467743796815756
$ws.Range("D3").Value = 2.048974129801823
$ws.Range("E3").Value = 1.620336905513141
$ws.Range("F3").Value = 1.296394889716211
$ws.Range("C4").Value = 2.320542303879344
$ws.Range("D4").Value = 1.95981595984354
$ws.Range("E4").Value = 1.42864171320842
$ws.Range("F4").Value = 1.258229492653228
$ws.Range("C5").Value = 2.230236732024127
$ws.Range("D5").Value = 1.871707064372008
$ws.Range("E5").Value = 1.303660677835469
$ws.Range("F5").Value = 1.212227701125545
$ws.Range("C6").Value = 2.152898553230203
$ws.Range("D6").Value = 1.83018001984209
$ws.Range("E6").Value = 1.23237768401677
$ws.Range("F6").Value = 1.18414203421364
$ws.Range("C7").Value = 2.112295704949706
$ws.Range("D7").Value = 1.777252250072877
$ws.Range("E7").Value = 1.211117753114425
$ws.Range("F7").Value = 1.176797248806555
$ws.Range("C8").Value = 2.958156851108843
$ws.Range("D8").Value = 2.674168078258705
$ws.Range("E8").Value = 2.550467843302049
$ws.Range("F8").Value = 1.583334841892314
$ws.Range("C9").Value = 2.752900061874218
$ws.Range("D9").Value = 2.583145584224332
$ws.Range("E9").Value = 2.273511081724153
$ws.Range("F9").Value = 1.493280646749689
$ws.Range("C10").Value = 2.582583693840313
$ws.Range("D10").Value = 2.475667281079886
$ws.Range("E10").Value = 2.100097495145042
$ws.Range("F10").Value = 1.427085768679523
$ws.Range("C11").Value = 2.448407076303334
$ws.Range("D11").Value = 2.367729447486429
$ws.Range("E11").Value = 1.93213609938342
$ws.Range("F11").Value = 1.365686344900394
$ws.Range("C12").Value = 2.372406859169245
$ws.Range("D12").Value = 2.309216619056827
$ws.Range("E12").Value = 1.811253008105346
$ws.Range("F12").Value = 1.354660165770392
$ws.Range("C13").Value = 2.350073057069686
$ws.Range("D13").Value = 2.292150511552243
$ws.Range("E13").Value = 1.78281293013583
$ws.Range("F13").Value = 1.345185568082276
$ws.Range("C14").Value = 4.186184925139913
$ws.Range("D14").Value = 3.705259749047888
$ws.Range("E14").Value = 3.997763106572068
$ws.Range("F14").Value = 7.516294411648735
$ws.Range("C15").Value = 3.529960235733572
$ws.Range("D15").Value = 3.219004613924047
$ws.Range("E15").Value = 3.393555781299532
$ws.Range("F15").Value = 4.42378924772177
$ws.Range("C16").Value = 3.047659905638771
$ws.Range("D16").Value = 3.022642601874604
$ws.Range("E16").Value = 2.860558958764877
$ws.Range("F16").Value = 4.818680357948059
$ws.Range("C17").Value = 2.737747394635913
$ws.Range("D17").Value = 2.896606764649546
$ws.Range("E17").Value = 2.63564813082195
$ws.Range("F17").Value = 1.554774843480185
$ws.Range("C18").Value = 2.604337145124109
$ws.Range("D18").Value = 2.841842788788624
$ws.Range("E18").Value = 2.508004185671087
$ws.Range("F18").Value = 5.621186927028381
$ws.Range("C19").Value = 2.560096222394556
$ws.Range("D19").Value = 2.830459825793664
$ws.Range("E19").Value = 2.482622735449927
$ws.Range("F19").Value = 5.122441358375069

$ws = $wb.Worksheets.Item("GWP")
$ws.Range("C2").Value = 3.177543609231863
$ws.Range("D2").Value = -0.03978436966721279
$ws.Range("E2").Value = -1.906796623613942
$ws.Range("F2").Value = 4.024350942332855
$ws.Range("C3").Value = 2.969266681453517
$ws.Range("D3").Value = -0.2311813185200753
$ws.Range("E3").Value = -2.204672901903208
$ws.Range("F3").Value = 3.898112040428388
$ws.Range("C4").Value = 2.797527053072463
$ws.Range("D4").Value = -0.3647913166204848
$ws.Range("E4").Value = -2.465470268185995
$ws.Range("F4").Value = 3.847423945034185
$ws.Range("C5").Value = 2.669637651768938
$ws.Range("D5").Value = -0.5114753847105374
$ws.Range("E5").Value = -2.634445552382329
$ws.Range("F5").Value = 3.790269888956434
$ws.Range("C6").Value = 2.585888885530413
$ws.Range("D6").Value = -0.647116785714429
$ws.Range("E6").Value = -2.751789613425808
$ws.Range("F6").Value = 3.784402992772982
$ws.Range("C7").Value = 2.562268440366928
$ws.Range("D7").Value = -0.6694730001569413
$ws.Range("E7").Value = -2.774240102855015
$ws.Range("F7").Value = 3.776778912376508
$ws.Range("C8").Value = 3.576287246536824
$ws.Range("D8").Value = 0.7101759992094349
$ws.Range("E8").Value = 0.4955796383988253
$ws.Range("F8").Value = 4.56119207966749
$ws.Range("C9").Value = 3.365879330218584
$ws.Range("D9").Value = 0.5930360215155734
$ws.Range("E9").Value = 0.3411104034775309
$ws.Range("F9").Value = 4.443628372338141
$ws.Range("C10").Value = 3.167043926109967
$ws.Range("D10").Value = 0.4626976785981737
$ws.Range("E10").Value = 0.2748553117789609
$ws.Range("F10").Value = 4.342497688574185
$ws.Range("C11").Value = 3.028202051342475
$ws.Range("D11").Value = 0.358776365922572
$ws.Range("E11").Value = 0.1225061720108877
$ws.Range("F11").Value = 4.252429318875691
$ws.Range("C12").Value = 2.951513286415557
$ws.Range("D12").Value = 0.2988633843486599
$ws.Range("E12").Value = 0.05446543762244752
$ws.Range("F12").Value = 4.257505494064553
$ws.Range("C13").Value = 2.928314674622027
$ws.Range("D13").Value = 0.2874365223886469
$ws.Range("E13").Value = 0.02043869269623629
$ws.Range("F13").Value = 4.247713267265833
$ws.Range("C14").Value = 4.209049280856689
$ws.Range("D14").Value = 1.750553115997926
$ws.Range("E14").Value = 2.330955706507889
$ws.Range("F14").Value = 8.025967677761347
$ws.Range("C15").Value = 3.953996322414121
$ws.Range("D15").Value = 1.761230797738574
$ws.Range("E15").Value = 2.237196927016431
$ws.Range("F15").Value = 7.540417283137299
$ws.Range("C16").Value = 3.763337137567906
$ws.Range("D16").Value = 1.607766186132912
$ws.Range("E16").Value = 2.151255212171631
$ws.Range("F16").Value = 5.819938901545838
$ws.Range("C17").Value = 3.618859103189898
$ws.Range("D17").Value = 1.560678251930118
$ws.Range("E17").Value = 2.142324663952231
$ws.Range("F17").Value = 4.867441673758091
$ws.Range("C18").Value = 3.527151274981034
$ws.Range("D18").Value = 1.554208836938743
$ws.Range("E18").Value = 2.062012077394965
$ws.Range("F18").Value = 5.774437242774585
$ws.Range("C19").Value = 3.501446753522099
$ws.Range("D19").Value = 1.538482330333837
$ws.Range("E19").Value = 2.049914574685602
$ws.Range("F19").Value = 5.738040569547072

$ws = $wb.Worksheets.Item("COD Price")
$ws.Range("C2").Value = 107.4296597082904
$ws.Range("D2").Value = -170.9490734061851
$ws.Range("E2").Value = -46.16986736214476
$ws.Range("F2").Value = 29.61220544483864
$ws.Range("C3").Value = 80.01418840355005
$ws.Range("D3").Value = -175.5594909888064
$ws.Range("E3").Value = -38.23076950325969
$ws.Range("F3").Value = -4.922154700028051
$ws.Range("C4").Value = 57.62675698839772
$ws.Range("D4").Value = -176.6592724885395
$ws.Range("E4").Value = -34.07260463176807
$ws.Range("F4").Value = -38.64690532393114
$ws.Range("C5").Value = 41.73340540086143
$ws.Range("D5").Value = -174.9432086932986
$ws.Range("E5").Value = -38.63969579411587
$ws.Range("F5").Value = -48.89885406738815
$ws.Range("C6").Value = 31.57732366730353
$ws.Range("D6").Value = -173.1412857688555
$ws.Range("E6").Value = -38.98735687509702
$ws.Range("F6").Value = -90.4157218435889
$ws.Range("C7").Value = 28.15463027314234
$ws.Range("D7").Value = -174.9348211391927
$ws.Range("E7").Value = -34.28704259331596
$ws.Range("F7").Value = -91.80752078820368
$ws.Range("C8").Value = 156.4924928251481
$ws.Range("D8").Value = 64.19398041408618
$ws.Range("E8").Value = 225.0057266317329
$ws.Range("F8").Value = 84.50975636875324
$ws.Range("C9").Value = 113.4578843523899
$ws.Range("D9").Value = 69.04200017046998
$ws.Range("E9").Value = 221.6319844641696
$ws.Range("F9").Value = 31.5568462593228
$ws.Range("C10").Value = 77.60461146999049
$ws.Range("D10").Value = 59.10119724568401
$ws.Range("E10").Value = 222.2795742563682
$ws.Range("F10").Value = -2.177439049266614
$ws.Range("C11").Value = 51.61839135927113
$ws.Range("D11").Value = 57.03780702584336
$ws.Range("E11").Value = 220.4671560861317
$ws.Range("F11").Value = -27.50425618612688
$ws.Range("C12").Value = 37.10988979800545
$ws.Range("D12").Value = 58.62836346204453
$ws.Range("E12").Value = 219.112031937348
$ws.Range("F12").Value = -45.17961087476454
$ws.Range("C13").Value = 32.36989815524446
$ws.Range("D13").Value = 56.80533239148262
$ws.Range("E13").Value = 220.1105575289303
$ws.Range("F13").Value = -51.40605853804824
$ws.Range("C14").Value = 438.00642728003
$ws.Range("D14").Value = 509.6693379311161
$ws.Range("E14").Value = 525.8922274528072
$ws.Range("F14").Value = 368.4414694810724
$ws.Range("C15").Value = 286.6411065117536
$ws.Range("D15").Value = 501.5524677971714
$ws.Range("E15").Value = 512.7388996244135
$ws.Range("F15").Value = 213.9424078540109
$ws.Range("C16").Value = 171.538325373381
$ws.Range("D16").Value = 511.5657045393625
$ws.Range("E16").Value = 520.774067300109
$ws.Range("F16").Value = 101.0910251704077
$ws.Range("C17").Value = 93.29437651703002
$ws.Range("D17").Value = 511.8725731928532
$ws.Range("E17").Value = 516.5485445287159
$ws.Range("F17").Value = 16.89823557386796
$ws.Range("C18").Value = 46.55291311604739
$ws.Range("D18").Value = 510.0274130729304
$ws.Range("E18").Value = 506.7975653304317
$ws.Range("F18").Value = -26.76540981788858
$ws.Range("C19").Value = 37.79185677047944
$ws.Range("D19").Value = 511.1445933820991
$ws.Range("E19").Value = 518.2404452495438
$ws.Range("F19").Value = -35.35592620018732

$ws = $wb.Worksheets.Item("COD GWP")
$ws.Range("C2").Value = -80.93300756433213
$ws.Range("D2").Value = -757.1966177259908
$ws.Range("E2").Value = -791.8671208545053
$ws.Range("F2").Value = -1481.164561885753
$ws.Range("C3").Value = -75.08284095990595
$ws.Range("D3").Value = -856.9401059815816
$ws.Range("E3").Value = -901.5787976304943
$ws.Range("F3").Value = -1548.126141576804
$ws.Range("C4").Value = -70.81979555539878
$ws.Range("D4").Value = -930.2206057029123
$ws.Range("E4").Value = -977.0006499052721
$ws.Range("F4").Value = -1605.316076231936
$ws.Range("C5").Value = -66.61113142272887
$ws.Range("D5").Value = -982.8920636381758
$ws.Range("E5").Value = -1036.581254894904
$ws.Range("F5").Value = -1646.414680078871
$ws.Range("C6").Value = -64.00145228219363
$ws.Range("D6").Value = -1030.894258805327
$ws.Range("E6").Value = -1074.999570201639
$ws.Range("F6").Value = -1673.399261691485
$ws.Range("C7").Value = -62.97892799691954
$ws.Range("D7").Value = -1030.899433085953
$ws.Range("E7").Value = -1085.049471906633
$ws.Range("F7").Value = -1681.792750930526
$ws.Range("C8").Value = -60.42308577715572
$ws.Range("D8").Value = -600.8513372340517
$ws.Range("E8").Value = -584.4185733315876
$ws.Range("F8").Value = -1233.466485286343
$ws.Range("C9").Value = -56.33059763461287
$ws.Range("D9").Value = -675.6420387036651
$ws.Range("E9").Value = -663.3919567223581
$ws.Range("F9").Value = -1301.148554351518
$ws.Range("C10").Value = -50.74692852206175
$ws.Range("D10").Value = -751.3006621524094
$ws.Range("E10").Value = -723.5853619289524
$ws.Range("F10").Value = -1350.13040345521
$ws.Range("C11").Value = -47.456425395276
$ws.Range("D11").Value = -798.4006048734512
$ws.Range("E11").Value = -767.4070385893466
$ws.Range("F11").Value = -1385.402829390831
$ws.Range("C12").Value = -45.15557654645145
$ws.Range("D12").Value = -828.5199526831867
$ws.Range("E12").Value = -804.7930228138418
$ws.Range("F12").Value = -1409.375447384455
$ws.Range("C13").Value = -44.1311259397657
$ws.Range("D13").Value = -837.2314875851117
$ws.Range("E13").Value = -809.652993471131
$ws.Range("F13").Value = -1415.893854833625
$ws.Range("C14").Value = -44.09613711799776
$ws.Range("D14").Value = -458.5911667135442
$ws.Range("E14").Value = -443.1269539172255
$ws.Range("F14").Value = -1038.219830213724
$ws.Range("C15").Value = -39.56980893090507
$ws.Range("D15").Value = -531.0718350507482
$ws.Range("E15").Value = -488.4288937746577
$ws.Range("F15").Value = -1089.154267887503
$ws.Range("C16").Value = -35.88459882019404
$ws.Range("D16").Value = -585.9988123076632
$ws.Range("E16").Value = -546.6070883231779
$ws.Range("F16").Value = -1129.359428272845
$ws.Range("C17").Value = -33.25388379961492
$ws.Range("D17").Value = -633.83093193697
$ws.Range("E17").Value = -592.6663543202593
$ws.Range("F17").Value = -1158.925164780158
$ws.Range("C18").Value = -31.35543874062113
$ws.Range("D18").Value = -645.5731067007855
$ws.Range("E18").Value = -615.3180842761971
$ws.Range("F18").Value = -1178.530209639965
$ws.Range("C19").Value = -30.65932278024449
$ws.Range("D19").Value = -668.3325194137763
$ws.Range("E19").Value = -622.9578686840193
$ws.Range("F19").Value = -1184.543716965397
